# Continue generate_practice_times_2: add the memb scheds to practice
# one by one (in the inner arrays) instead of a separate method.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet5")

$xlPasteFormats = -4122

# Add "Quynh-Nhi" row first (ends up as row 9 once Karen's row is
# inserted above it).
$ws.Range("A9").Value = "Quynh-Nhi"
$ws.Range("B9").Value = "Free"
$ws.Range("C9").Value = "After 6"
$ws.Range("D9").Value = "After 6"
$ws.Range("E9").Value = "After 6"
$ws.Range("F9").Value = "After 6"
$ws.Range("G9").Value = "After 4"
$ws.Range("H9").Value = "Free"

# Match row 9 formatting to the existing schedule rows.
$ws.Range("C7:H7").Copy()
$ws.Range("A9").PasteSpecial($xlPasteFormats)
$ws.Range("C7:H7").Copy()
$ws.Range("B9:H9").PasteSpecial($xlPasteFormats)

# Row 8 previously held "michelle f" - replace it entirely with
# "Karen"'s schedule (no trailing OTHER/comment column for her).
$ws.Range("A8").Value = "Karen"
$ws.Range("B8").Value = "After 6pm"
$ws.Range("C8").Value = "After 3:30"
$ws.Range("D8").Value = "Not Available"
$ws.Range("E8").Value = "Not available"
$ws.Range("F8").Value = "After 3:30"
$ws.Range("G8").Value = "After 3:30"
$ws.Range("H8").Value = "After 5:30"
$ws.Range("I8").Clear()

# Normalize D8:G8 formatting to match the rest of the row (was inherited
# from the old "michelle f" row's distinct style).
$ws.Range("C8").Copy()
$ws.Range("D8:G8").PasteSpecial($xlPasteFormats)

# Row 6 (Cindy): update Saturday then Wednesday availability text.
$ws.Range("H6").Value = "Free except 4h30-6h30"
$ws.Range("E6").Value = "free after 4h30"

$ws.Range("E7").Select()
